$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Set-RowValues {
    param($ws, $rowNum, $values)
    $arr = New-Object "object[,]" 1,($values.Count)
    for ($i = 0; $i -lt $values.Count; $i++) {
        $arr[0,$i] = $values[$i]
    }
    $ws.Range("A$rowNum" + ":AC$rowNum").Value2 = $arr
}

# New rows 115-118 do not yet exist in the sheet; clone the cell
# formatting (id column + date column) from an existing data row first,
# so the new rows look consistent with the rest of the table.
$ws.Range("A4").Copy() | Out-Null
$ws.Range("A115:A118").PasteSpecial(-4122) | Out-Null
$ws.Range("E4").Copy() | Out-Null
$ws.Range("E115:E118").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

Set-RowValues $ws 4 @(2, 6139017, 'Estonia Meistriliiga', 'Estonia Meistriliiga', 45084.5, 'JK Tammeka Tartu', 'Harju JK Laagri', 2, 0, 'H', 1.666, 3.6, 4.2, 1.727, 3.5, 4, -0.75, 2, 1.8, 2.5, 1.9, 1.9, 0.7270000000000001, -1, -1, 1, -1, -1, 0.8999999999999999)
Set-RowValues $ws 5 @(3, 6139018, 'Estonia Meistriliiga', 'Estonia Meistriliiga', 45084.5, 'JK Tallinna Kalev', 'JK Trans Narva', 0, 1, 'A', 2.4, 3.4, 2.5, 2.875, 3.1, 2.3, 0.25, 1.75, 2.05, 2.25, 1.925, 1.875, -1, -1, 1.3, -1, 1.05, -1, 0.875)
Set-RowValues $ws 10 @(8, 6139020, 'Estonia Meistriliiga', 'Estonia Meistriliiga', 45088.45833333334, 'FC Levadia Tallinn', 'JK Tallinna Kalev', 2, 1, 'H', 1.363, 4.333, 6.5, 1.3, 4.75, 8, -1.75, 1.95, 1.85, 2.75, 1.8, 2, 0.3, -1, -1, -1, 0.8500000000000001, 0.4, -0.5)
Set-RowValues $ws 11 @(9, 6139023, 'Estonia Meistriliiga', 'Estonia Meistriliiga', 45088.54166666666, 'JK Trans Narva', 'JK Tammeka Tartu', 2, 0, 'H', 2, 3.3, 3.2, 2.05, 3.3, 3.1, -0.25, 1.85, 1.95, 2.25, 1.925, 1.875, 1.05, -1, -1, 0.8500000000000001, -1, -0.5, 0.4375)
Set-RowValues $ws 13 @(11, 6139025, 'Estonia Meistriliiga', 'Estonia Meistriliiga', 45104.54166666666, 'JK Tammeka Tartu', 'Parnu JK Vaprus', 2, 3, 'A', 2.2, 3.3, 2.8, 1.7, 3.6, 4.333, -0.75, 1.95, 1.85, 2.25, 1.875, 1.925, -1, -1, 3.333, -1, 0.8500000000000001, 0.875, -1)
Set-RowValues $ws 15 @(13, 6139026, 'Estonia Meistriliiga', 'Estonia Meistriliiga', 45105.5, 'JK Tallinna Kalev', 'Harju JK Laagri', 1, 1, 'D', 1.444, 4.75, 5, 1.666, 4.5, 3.6, -0.75, 1.875, 1.925, 2.5, 1.875, 1.925, -1, 3.5, -1, -1, 0.925, -1, 0.925)
Set-RowValues $ws 20 @(18, 6138125, 'Estonia Meistriliiga', 'Estonia Meistriliiga', 45109.45833333334, 'JK Nomme Kalju', 'JK Tallinna Kalev', 1, 2, 'A', 1.571, 3.8, 4.75, 1.55, 3.8, 5, -1, 1.95, 1.85, 2.75, 1.95, 1.85, -1, -1, 4, -1, 0.8500000000000001, 0.475, -0.5)
Set-RowValues $ws 21 @(19, 6139031, 'Estonia Meistriliiga', 'Estonia Meistriliiga', 45109.54166666666, 'FC Flora Tallinn', 'JK Tammeka Tartu', 3, 0, 'H', 1.2, 6, 9, 1.125, 8, 17, -2.5, 1.925, 1.875, 3.5, 1.975, 1.825, 0.125, -1, -1, 0.925, -1, -1, 0.825)
Set-RowValues $ws 24 @(22, 6139034, 'Estonia Meistriliiga', 'Estonia Meistriliiga', 45115.5625, 'JK Tallinna Kalev', 'FC Kuressaare', 1, 3, 'A', 2.75, 3.5, 2.2, 3.1, 3.5, 2.05, 0.5, 1.725, 2.075, 2.75, 1.925, 1.875, -1, -1, 1.05, -1, 1.075, 0.925, -1)
Set-RowValues $ws 25 @(23, 6139035, 'Estonia Meistriliiga', 'Estonia Meistriliiga', 45116.45833333334, 'Harju JK Laagri', 'JK Tammeka Tartu', 3, 0, 'H', 3, 3.5, 2.05, 2.55, 3.5, 2.375, 0, 1.975, 1.825, 2.5, 1.9, 1.9, 1.55, -1, -1, 0.9750000000000001, -1, 0.8999999999999999, -1)
Set-RowValues $ws 27 @(25, 6139037, 'Estonia Meistriliiga', 'Estonia Meistriliiga', 45129.45833333334, 'JK Tallinna Kalev', 'Parnu JK Vaprus', 1, 1, 'D', 2.2, 3.3, 2.8, 2.1, 3.3, 3, -0.25, 1.875, 1.925, 2.5, 1.875, 1.925, -1, 2.3, -1, -0.5, 0.4625, -1, 0.925)
Set-RowValues $ws 30 @(28, 6139039, 'Estonia Meistriliiga', 'Estonia Meistriliiga', 45130.5, 'JK Tammeka Tartu', 'FC Levadia Tallinn', 0, 0, 'D', 6.5, 5.5, 1.285, 8.5, 6.5, 1.2, 2, 1.8, 2, 3, 2.025, 1.775, -1, 5.5, -1, 0.8, -1, -1, 0.7749999999999999)
Set-RowValues $ws 33 @(31, 6825464, 'Estonia Meistriliiga', 'Estonia Meistriliiga', 45138.57291666666, 'FC Kuressaare', 'JK Tammeka Tartu', 3, 1, 'H', 1.85, 3.6, 3.4, 1.95, 3.5, 3.25, -0.5, 2, 1.8, 2.5, 1.975, 1.825, 0.95, -1, -1, 1, -1, 0.9750000000000001, -1)
Set-RowValues $ws 36 @(34, 6139041, 'Estonia Meistriliiga', 'Estonia Meistriliiga', 45143.5625, 'JK Tallinna Kalev', 'JK Tammeka Tartu', 1, 1, 'D', 1.8, 3.75, 3.4, 2.1, 3.5, 2.8, -0.25, 1.9, 1.9, 2.5, 1.9, 1.9, -1, 2.5, -1, -0.5, 0.45, -1, 0.8999999999999999)
Set-RowValues $ws 39 @(37, 6139045, 'Estonia Meistriliiga', 'Estonia Meistriliiga', 45150.45833333334, 'Paide Linnameeskond', 'JK Tallinna Kalev', 1, 1, 'D', 1.533, 3.75, 5, 1.444, 4.2, 6, -1.25, 2, 1.8, 2.75, 1.95, 1.85, -1, 3.2, -1, -1, 0.8, -1, 0.8500000000000001)
Set-RowValues $ws 41 @(39, 6138128, 'Estonia Meistriliiga', 'Estonia Meistriliiga', 45151.45833333334, 'JK Tammeka Tartu', 'JK Nomme Kalju', 1, 1, 'D', 2.75, 3.2, 2.3, 3.75, 3.4, 1.8, 0.5, 1.95, 1.85, 2.25, 1.75, 1.95, -1, 2.4, -1, 0.95, -1, -0.5, 0.475)
Set-RowValues $ws 44 @(42, 6139050, 'Estonia Meistriliiga', 'Estonia Meistriliiga', 45157.45833333334, 'Paide Linnameeskond', 'JK Tammeka Tartu', 6, 3, 'H', 1.45, 4, 6.5, 1.5, 4, 5.75, -1, 1.875, 1.925, 2.5, 2.025, 1.775, 0.5, -1, -1, 0.875, -1, 1.025, -1)
Set-RowValues $ws 47 @(45, 6139051, 'Estonia Meistriliiga', 'Estonia Meistriliiga', 45158.5, 'JK Tallinna Kalev', 'FC Flora Tallinn', 1, 1, 'D', 10, 6, 1.2, 10, 5.5, 1.222, 1.75, 1.975, 1.825, 3, 1.95, 1.85, -1, 4.5, -1, 0.9750000000000001, -1, -1, 0.8500000000000001)
Set-RowValues $ws 49 @(47, 6139053, 'Estonia Meistriliiga', 'Estonia Meistriliiga', 45164.45833333334, 'JK Tammeka Tartu', 'Parnu JK Vaprus', 0, 0, 'D', 2.5, 3.4, 2.375, 2.05, 3.4, 3, -0.25, 1.825, 1.975, 2.5, 1.975, 1.825, -1, 2.4, -1, -0.5, 0.4875, -1, 0.825)
Set-RowValues $ws 50 @(48, 6139054, 'Estonia Meistriliiga', 'Estonia Meistriliiga', 45164.54166666666, 'FC Levadia Tallinn', 'JK Tallinna Kalev', 1, 1, 'D', 1.3, 4.333, 9, 1.285, 4.75, 9, -1.75, 1.975, 1.825, 3, 1.925, 1.875, -1, 3.75, -1, -1, 0.825, -1, 0.875)
Set-RowValues $ws 53 @(51, 6139056, 'Estonia Meistriliiga', 'Estonia Meistriliiga', 45170.5, 'JK Trans Narva', 'JK Tallinna Kalev', 2, 1, 'H', 2.375, 3.4, 2.625, 2.2, 3.5, 2.8, -0.25, 1.975, 1.825, 2.5, 1.825, 1.975, 1.2, -1, -1, 0.9750000000000001, -1, 0.825, -1)
Set-RowValues $ws 55 @(53, 6139057, 'Estonia Meistriliiga', 'Estonia Meistriliiga', 45171.45833333334, 'FC Kuressaare', 'JK Tammeka Tartu', 1, 1, 'D', 1.833, 3.5, 4, 2.05, 3.4, 3.25, -0.25, 1.8, 2, 2.75, 1.8, 2, -1, 2.4, -1, -0.5, 0.5, -1, 1)
Set-RowValues $ws 58 @(56, 6139060, 'Estonia Meistriliiga', 'Estonia Meistriliiga', 45184.54166666666, 'JK Tammeka Tartu', 'JK Trans Narva', 3, 0, 'H', 2.75, 3.25, 2.375, 2.7, 3.3, 2.375, 0, 2.025, 1.775, 2.5, 1.875, 1.925, 1.7, -1, -1, 1.025, -1, 0.875, -1)
Set-RowValues $ws 59 @(57, 6139061, 'Estonia Meistriliiga', 'Estonia Meistriliiga', 45185.45833333334, 'JK Tallinna Kalev', 'Harju JK Laagri', 2, 1, 'H', 1.727, 3.4, 4.6, 1.7, 3.5, 4.75, -0.75, 1.85, 1.95, 2.5, 1.85, 1.95, 0.7, -1, -1, 0.425, -0.5, 0.8500000000000001, -1)
Set-RowValues $ws 63 @(61, 6139068, 'Estonia Meistriliiga', 'Estonia Meistriliiga', 45191.57291666666, 'FC Flora Tallinn', 'JK Tallinna Kalev', 1, 0, 'H', 1.166, 7, 10, 1.142, 7, 12, -2.25, 1.975, 1.825, 3, 1.8, 2, 0.1419999999999999, -1, -1, -1, 0.825, -1, 1)
Set-RowValues $ws 64 @(62, 6139067, 'Estonia Meistriliiga', 'Estonia Meistriliiga', 45192.35416666666, 'Paide Linnameeskond', 'Parnu JK Vaprus', 3, 0, 'H', 1.8, 3.4, 3.8, 1.5, 3.8, 6, -1, 1.75, 1.95, 2.5, 1.9, 1.9, 0.5, -1, -1, 0.75, -1, 0.8999999999999999, -1)
Set-RowValues $ws 65 @(63, 6139064, 'Estonia Meistriliiga', 'Estonia Meistriliiga', 45192.35416666666, 'JK Trans Narva', 'Harju JK Laagri', 1, 3, 'A', 1.75, 3.6, 3.8, 1.45, 4, 6, -1, 1.85, 1.95, 2.5, 1.9, 1.9, -1, -1, 5, -1, 0.95, 0.8999999999999999, -1)
Set-RowValues $ws 66 @(64, 6138133, 'Estonia Meistriliiga', 'Estonia Meistriliiga', 45192.45833333334, 'JK Nomme Kalju', 'JK Tammeka Tartu', 4, 1, 'H', 1.75, 3.4, 4, 1.533, 3.5, 5.25, -1, 2, 1.8, 2.5, 1.85, 1.95, 0.5329999999999999, -1, -1, 1, -1, 0.8500000000000001, -1)
Set-RowValues $ws 70 @(68, 6139070, 'Estonia Meistriliiga', 'Estonia Meistriliiga', 45196.45833333334, 'JK Tallinna Kalev', 'Paide Linnameeskond', 2, 1, 'H', 5, 4, 1.5, 3.75, 3.75, 1.75, 0.75, 1.825, 1.975, 2.75, 2, 1.8, 2.75, -1, -1, 0.825, -1, 0.5, -0.5)
Set-RowValues $ws 72 @(70, 6139072, 'Estonia Meistriliiga', 'Estonia Meistriliiga', 45196.54166666666, 'JK Tammeka Tartu', 'FC Flora Tallinn', 1, 2, 'A', 9, 7, 1.166, 7, 6, 1.25, 1.75, 1.9, 1.9, 3, 1.95, 1.85, -1, -1, 0.25, 0.8999999999999999, -1, 0, 0)
Set-RowValues $ws 74 @(72, 6139073, 'Estonia Meistriliiga', 'Estonia Meistriliiga', 45199.45833333334, 'JK Tammeka Tartu', 'Paide Linnameeskond', 1, 2, 'A', 5, 4, 1.5, 4.5, 3.8, 1.571, 0.75, 1.95, 1.75, 2.75, 1.9, 1.9, -1, -1, 0.571, -0.5, 0.375, 0.45, -0.5)
Set-RowValues $ws 75 @(73, 6139075, 'Estonia Meistriliiga', 'Estonia Meistriliiga', 45200.35416666666, 'FC Kuressaare', 'JK Tallinna Kalev', 1, 1, 'D', 3.6, 3.6, 1.8, 2.625, 3.5, 2.3, 0.25, 1.725, 1.975, 2.75, 1.8, 2, -1, 2.5, -1, 0.3625, -0.5, -1, 1)
Set-RowValues $ws 79 @(77, 6354607, 'Estonia Meistriliiga', 'Estonia Meistriliiga', 45206.35416666666, 'JK Tallinna Kalev', 'FC Levadia Tallinn', 1, 2, 'A', 6.5, 5.5, 1.285, 6.5, 5.5, 1.285, 1.5, 1.95, 1.85, 3, 1.95, 1.75, -1, -1, 0.2849999999999999, 0.95, -1, 0, 0)
Set-RowValues $ws 80 @(78, 6368429, 'Estonia Meistriliiga', 'Estonia Meistriliiga', 45206.45833333334, 'JK Trans Narva', 'JK Tammeka Tartu', 2, 0, 'H', 2, 3.3, 3.2, 2.15, 3.3, 2.9, -0.25, 1.925, 1.875, 2.5, 1.85, 1.95, 1.15, -1, -1, 0.925, -1, -1, 0.95)
Set-RowValues $ws 83 @(81, 6368431, 'Estonia Meistriliiga', 'Estonia Meistriliiga', 45220.27083333334, 'JK Tallinna Kalev', 'Parnu JK Vaprus', 1, 0, 'H', 2.1, 3.4, 2.9, 1.95, 3.4, 3.4, -0.5, 2, 1.8, 2.5, 1.925, 1.875, 0.95, -1, -1, 1, -1, -1, 0.875)
Set-RowValues $ws 85 @(83, 6369469, 'Estonia Meistriliiga', 'Estonia Meistriliiga', 45220.45833333334, 'FC Levadia Tallinn', 'JK Tammeka Tartu', 2, 1, 'H', 1.2, 6, 9, 1.222, 5.5, 10, -1.75, 1.8, 2, 3, 1.9, 1.9, 0.222, -1, -1, -1, 1, 0, 0)
Set-RowValues $ws 88 @(86, 6376947, 'Estonia Meistriliiga', 'Estonia Meistriliiga', 45226.54166666666, 'JK Tammeka Tartu', 'JK Tallinna Kalev', 2, 7, 'A', 3.6, 3.4, 1.909, 2.4, 3.6, 2.45, 0, 1.875, 1.925, 2.75, 1.975, 1.825, -1, -1, 1.45, -1, 0.925, 0.9750000000000001, -1)
Set-RowValues $ws 93 @(91, 6418047, 'Estonia Meistriliiga', 'Estonia Meistriliiga', 45230.54166666666, 'Harju JK Laagri', 'JK Tallinna Kalev', 0, 2, 'A', 3.75, 3.75, 1.727, 2.875, 3.3, 2.15, 0.25, 1.85, 1.95, 2.5, 1.875, 1.925, -1, -1, 1.15, -1, 0.95, -1, 0.925)
Set-RowValues $ws 96 @(94, 6482819, 'Estonia Meistriliiga', 'Estonia Meistriliiga', 45231.54166666666, 'JK Tammeka Tartu', 'FC Kuressaare', 0, 1, 'A', 1.833, 3.5, 3.5, 2.1, 3.4, 2.875, -0.25, 1.975, 1.825, 3, 1.825, 1.975, -1, -1, 1.875, -1, 0.825, -1, 0.9750000000000001)
Set-RowValues $ws 100 @(98, 6528945, 'Estonia Meistriliiga', 'Estonia Meistriliiga', 45235.3125, 'JK Tallinna Kalev', 'JK Nomme Kalju', 1, 1, 'D', 4, 4, 1.615, 3.4, 3.75, 1.8, 0.5, 1.975, 1.825, 2.75, 1.8, 2, -1, 2.75, -1, 0.9750000000000001, -1, -1, 1)
Set-RowValues $ws 102 @(100, 6533011, 'Estonia Meistriliiga', 'Estonia Meistriliiga', 45235.5, 'Parnu JK Vaprus', 'JK Tammeka Tartu', 2, 0, 'H', 2.3, 3.4, 2.625, 2.3, 3.4, 2.625, -0.25, 1.975, 1.725, 2.5, 1.775, 1.925, 1.3, -1, -1, 0.9750000000000001, -1, -1, 0.925)
Set-RowValues $ws 103 @(101, 7440206, 'Estonia Meistriliiga', 'Estonia Meistriliiga', 45240.58333333334, 'JK Tammeka Tartu', 'Harju JK Laagri', 2, 1, 'H', 1.909, 4, 3.1, 1.95, 3.8, 3, -0.25, 1.8, 2, 2.75, 1.975, 1.825, 0.95, -1, -1, 0.8, -1, 0.4875, -0.5)
Set-RowValues $ws 105 @(103, 6537869, 'Estonia Meistriliiga', 'Estonia Meistriliiga', 45241.375, 'JK Tallinna Kalev', 'JK Trans Narva', 5, 0, 'H', 1.6, 4, 4.5, 1.65, 4, 4.333, -0.75, 1.8, 2, 2.75, 1.9, 1.9, 0.6499999999999999, -1, -1, 0.8, -1, 0.8999999999999999, -1)
Set-RowValues $ws 106 @(104, 6535416, 'Estonia Meistriliiga', 'Estonia Meistriliiga', 45241.375, 'Paide Linnameeskond', 'FC Levadia Tallinn', 2, 2, 'D', 3, 3.8, 2, 3, 4, 1.909, 0.5, 1.85, 1.95, 2.75, 1.95, 1.85, -1, 3, -1, 0.8500000000000001, -1, 0.95, -1)
Set-RowValues $ws 107 @(105, 6537957, 'Estonia Meistriliiga', 'Estonia Meistriliiga', 45241.375, 'FC Flora Tallinn', 'JK Nomme Kalju', 0, 0, 'D', 1.4, 4, 7.5, 1.5, 4.2, 5, -1, 1.85, 1.95, 2.75, 1.85, 1.95, -1, 3.2, -1, -1, 0.95, -1, 0.95)
Set-RowValues $ws 108 @(106, 7551820, 'Estonia Meistriliiga', 'Estonia Meistriliiga', 45263.54166666666, 'JK Tammeka Tartu', 'Viimsi JK', 1, 1, 'D', 1.285, 5, 7.5, 1.285, 4.75, 8, -1.75, 2, 1.8, 3, 1.925, 1.875, -1, 3.75, -1, -1, 0.8, -1, 0.875)
Set-RowValues $ws 112 @(110, 7891675, 'Estonia Meistriliiga', 'Estonia Meistriliiga', 45354.3125, 'JK Tallinna Kalev', 'JK Tammeka Tartu', 1, 1, 'D', 1.571, 3.75, 4.75, 1.65, 3.8, 4.333, -0.75, 1.9, 1.9, 2.75, 1.875, 1.925, -1, 2.8, -1, -1, 0.8999999999999999, -1, 0.925)
Set-RowValues $ws 114 @(112, 7919321, 'Estonia Meistriliiga', 'Estonia Meistriliiga', 45360.3125, 'JK Tammeka Tartu', 'Parnu JK Vaprus', 1, 0, 'H', 2.4, 3.6, 2.4, 2.1, 3.6, 2.8, -0.25, 1.9, 1.9, 2.5, 1.9, 1.9, 1.1, -1, -1, 0.8999999999999999, -1, -1, 0.8999999999999999)
Set-RowValues $ws 115 @(113, 7919323, 'Estonia Meistriliiga', 'Estonia Meistriliiga', 45360.39583333334, 'JK Nomme Kalju', 'JK Trans Narva', 3, 0, 'H', 1.285, 5.5, 6.5, 1.571, 4.75, 4.2, -1, 1.925, 1.875, 2.75, 1.875, 1.925, 0.571, -1, -1, 0.925, -1, 0.4375, -0.5)
Set-RowValues $ws 116 @(114, 7919322, 'Estonia Meistriliiga', 'Estonia Meistriliiga', 45360.39583333334, 'FC Kuressaare', 'FC Levadia Tallinn', 0, 6, 'A', 11, 6, 1.166, 15, 8.5, 1.125, 2.5, 1.825, 1.975, 3.25, 1.9, 1.9, -1, -1, 0.125, -1, 0.9750000000000001, 0.8999999999999999, -1)
Set-RowValues $ws 117 @(115, 7721006, 'Estonia Meistriliiga', 'Estonia Meistriliiga', 45361.3125, 'Paide Linnameeskond', 'JK Nomme United', 3, 1, 'H', 1.285, 5.5, 6.5, 1.181, 6.5, 9, -2, 1.825, 1.975, 3.5, 1.95, 1.85, 0.181, -1, -1, 0, 0, 0.95, -1)
Set-RowValues $ws 118 @(116, 7719642, 'Estonia Meistriliiga', 'Estonia Meistriliiga', 45361.39583333334, 'FC Flora Tallinn', 'JK Tallinna Kalev', 2, 2, 'D', 1.444, 4, 6, 1.444, 3.8, 6, -1.25, 1.975, 1.825, 2.75, 1.9, 1.9, -1, 2.8, -1, -1, 0.825, 0.8999999999999999, -1)

$wb.Save()
